# Trading update: 2026-02-18 10:38:03
#
# 1) Drop the "momentum" strategy sheet; the live per-strategy tab position
#    it used to occupy is reused for the (now sole) "MarketMaking" sheet.
# 2) The old, separate "MarketMaking" sheet (5 open trades) is replaced by a
#    fresh single-open-trade snapshot (trade #31).
# 3) "All Trades" gets its trailing header columns reordered (Capital
#    After is now followed by Entry Slippage / Exit Slippage / Confidence /
#    Entry Reason / Exit Reason / Duration), every still-OPEN trade's
#    now-stale derived columns (Capital After .. Duration) are cleared back
#    out, and the new trade #31 is appended as row 32.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: remove the old, separate "MarketMaking" sheet entirely (it will
# be reconstituted from the "momentum" sheet below with fresh data).
# ---------------------------------------------------------------------
$oldMarketMaking = $wb.Worksheets.Item("MarketMaking")
$oldMarketMaking.Delete()

# ---------------------------------------------------------------------
# Step 2: turn the "momentum" sheet into the new "MarketMaking" sheet:
# drop its 3 stale data rows and replace them with the single new trade.
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("momentum")

$mm.Rows.Item(3).EntireRow.Delete()
$mm.Rows.Item(3).EntireRow.Delete()

$mm.Range("A2").Value = 31
$mm.Range("B2").NumberFormat = "@"
$mm.Range("B2").Value = "2026-02-18"
$mm.Range("C2").NumberFormat = "@"
$mm.Range("C2").Value = "10:37:52"
$mm.Range("D2").Value = "MarketMaking"
$mm.Range("E2").Value = "UP"
$mm.Range("F2").Value = 0.6899999999999999
$mm.Range("G2").Value = ""
$mm.Range("H2").Value = "OPEN"
$mm.Range("I2").Value = 0
$mm.Range("J2").Value = 0
$mm.Range("K2").Value = 100
$mm.Range("L2").Value = 0
$mm.Range("M2").Value = 0
$mm.Range("N2").Value = 0.6
$mm.Range("O2").Value = "Normal spread capture: 194 bps"
$mm.Range("P2").Value = ""
$mm.Range("Q2").Value = 0

$mm.Name = "MarketMaking"

# ---------------------------------------------------------------------
# Step 3: "All Trades" sheet updates.
# ---------------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

# 3a. Reorder the trailing headers (L1:Q1).
$all.Range("L1").Value = "Entry Slippage (bps)"
$all.Range("M1").Value = "Exit Slippage (bps)"
$all.Range("N1").Value = "Confidence"
$all.Range("O1").Value = "Entry Reason"
$all.Range("P1").Value = "Exit Reason"
$all.Range("Q1").Value = "Duration (min)"

# 3b. Rows 6 & 7 (trades 5 & 6): clear the now-reassigned K:M columns.
foreach ($r in 6,7) {
    $all.Range("K$r").Value = ""
    $all.Range("L$r").Value = ""
    $all.Range("M$r").Value = ""
}

# 3c. Rows 24-31 (trades 23-30, still OPEN): Exit Price becomes an
# explicit 0 and every derived column K:Q is cleared back to blank.
foreach ($r in 24,25,26,27,28,29,30,31) {
    $all.Range("G$r").Value = 0
    foreach ($col in "K","L","M","N","O","P","Q") {
        $all.Range("$col$r").Value = ""
    }
}

# 3d. Append the new trade (#31) as row 32.
$all.Range("A32").Value = 31
$all.Range("B32").NumberFormat = "@"
$all.Range("B32").Value = "2026-02-18"
$all.Range("C32").NumberFormat = "@"
$all.Range("C32").Value = "10:37:52"
$all.Range("D32").Value = "MarketMaking"
$all.Range("E32").Value = "UP"
$all.Range("F32").Value = 0.6899999999999999
$all.Range("G32").Value = ""
$all.Range("H32").Value = "OPEN"
$all.Range("I32").Value = 0
$all.Range("J32").Value = 0
$all.Range("K32").Value = 100
$all.Range("L32").Value = 0
$all.Range("M32").Value = 0
$all.Range("N32").Value = 0.6
$all.Range("O32").Value = "Normal spread capture: 194 bps"
$all.Range("P32").Value = ""
$all.Range("Q32").Value = 0
